$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the next backlog item as row 14, cloning the formatting of the
# previous data row (grey font, date number format) before filling values.
$ws.Range("A11:E11").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A14").Value = 13
$ws.Range("B14").Value2 = 42708
$ws.Range("C14").Value = "Fix Test Scripts bug with echoing test totals"
$ws.Range("E14").Value = "COMPLETED"

# Leave the selection where the author left it after entering the row
$ws.Range("D11").Select()
